# Fixed a conversion fuck up
#
# The workbook contained formulas/defined names that referenced an external
# workbook (RCL-B-COM1.xlsx). That external reference had gone stale and
# every formula depending on it was showing #REF!. This script breaks the
# now-dead external link and replaces the broken formulas with the literal
# values they used to resolve to, then lets the normal formulas that were
# merely *downstream* of those broken cells recompute cleanly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the (now unused) external workbook link, which also drops the
# <externalReferences> block from workbook.xml and the externalLink1.xml part.
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# C7: float precision clean-up (0.30000000000000004 -> 0.3)
$ws.Range("C7").Value = 0.3

# C8: was "=loss" (#REF!) -> literal value
$ws.Range("C8").Formula = "0.2"

# C12: was "=C4-C5-10*LOG(SSNT)" (#REF!) -> SSNT replaced with literal 261
$ws.Range("C12").Formula = "=C4-C5-10*LOG(261)"

# C19: was "=SCG" (#REF!) -> literal value
$ws.Range("C19").Formula = "5"

# C20: was "=SCL" (#REF!) -> literal value
$ws.Range("C20").Formula = "2"

# C22 / C23: updated constants
$ws.Range("C22").Value = 6.8
$ws.Range("C23").Value = 0.1

# C27: was "=C19-C20-10*LOG(SSNT)" (#REF!) -> SSNT replaced with literal 261
$ws.Range("C27").Formula = "=C19-C20-10*LOG(261)"

# The cells below never referenced the broken names directly, but their
# cached values were #REF! because they depend (directly or transitively)
# on the cells fixed above. Re-apply their existing formulas so they get
# re-evaluated against the corrected inputs.
$ws.Range("C11").Formula = "=C6-SUM(C7:C8,C10)"
$ws.Range("C13").Formula = "=C11-C7--228.6+C12"
$ws.Range("C15").Formula = "=C13-C14"
$ws.Range("C16").Formula = "=0.5*ERFC(2*(C15/SQRT(2)))"
$ws.Range("C21").Formula = "=C18+C19-C20"
$ws.Range("C26").Formula = "=C21-SUM(C22:C23,C25)"
$ws.Range("C28").Formula = "=C26-C22--228.6+C27"
$ws.Range("C30").Formula = "=C28-C29"
$ws.Range("C31").Formula = "=0.5*ERFC(2*(C30/SQRT(2)))"

$excel.CalculateFull()

# Restore the last active selection on Sheet1
$ws.Range("C24").Select()

$wb.Save()
